# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Hoja1: actualizar el resumen de conversión del día (celda A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$old1 = "1000 Bs = 14.36 = 59052.4 pesos"
$new1 = "1000 Bs = 14.15 = 57973.25 pesos"
$old2 = "59052.4 pesos = 14.32 = 973.71 Bs"
$new2 = "57973.25 pesos = 14.07 = 972.57 Bs"

$texto = $wsHoja1.Range("A1").Value2
$texto = $texto.Replace($old1, $new1)
$texto = $texto.Replace($old2, $new2)
$wsHoja1.Range("A1").Value = $texto

# --- tasas: actualizar tasas Binance / transfi (N10, O10, N12, O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 70.65000000000001
$wsTasas.Range("O10").Value = 4095.81
$wsTasas.Range("N12").Value = 4119
$wsTasas.Range("O12").Value = 69.101
